$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("B1").Value = "二氧化硅(SiO2)"
$ws.Range("C1").Value = "氧化钠(Na2O)"
$ws.Range("D1").Value = "氧化钾(K2O)"
$ws.Range("E1").Value = "氧化钙(CaO)"
$ws.Range("F1").Value = "氧化镁(MgO)"
$ws.Range("G1").Value = "氧化铝(Al2O3)"
$ws.Range("H1").Value = "氧化铁(Fe2O3)"
$ws.Range("I1").Value = "氧化铜(CuO)"
$ws.Range("J1").Value = "氧化铅(PbO)"
$ws.Range("K1").Value = "氧化钡(BaO)"
$ws.Range("L1").Value = "五氧化二磷(P2O5)"
$ws.Range("M1").Value = "氧化锶(SrO)"
$ws.Range("N1").Value = "氧化锡(SnO2)"
$ws.Range("O1").Value = "二氧化硫(SO2)"

# Update data rows 2-7 with new computed values
$ws.Range("B2").Value = 67.75905951623973
$ws.Range("C2").Value = 1.162755650491632
$ws.Range("D2").Value = 10.28666759733919
$ws.Range("E2").Value = 5.541647314699883
$ws.Range("F2").Value = 1.091688068066482
$ws.Range("G2").Value = 6.729608836911932
$ws.Range("H2").Value = 1.929707610322754
$ws.Range("I2").Value = 2.450747743951087
$ws.Range("J2").Value = 0.5451869611956746
$ws.Range("K2").Value = 0.818801560937361
$ws.Range("L2").Value = 1.515321478859036
$ws.Range("M2").Value = 0.05228361519063392
$ws.Range("N2").Value = 0.3945951823073121
$ws.Range("O2").Value = 0.3249974193755092
$ws.Range("B3").Value = 67.77407670251333
$ws.Range("C3").Value = 1.151290898922405
$ws.Range("D3").Value = 10.27580511324248
$ws.Range("E3").Value = 5.526944620243655
$ws.Range("F3").Value = 1.089097736936
$ws.Range("G3").Value = 6.736879365978024
$ws.Range("H3").Value = 1.898532392881102
$ws.Range("I3").Value = 2.451019082829195
$ws.Range("J3").Value = 0.5553877199370215
$ws.Range("K3").Value = 0.793416431564602
$ws.Range("L3").Value = 1.511979004014555
$ws.Range("M3").Value = 0.05236598373623543
$ws.Range("N3").Value = 0.4235211913896434
$ws.Range("O3").Value = 0.3297057604488701
$ws.Range("B4").Value = 67.78398999069908
$ws.Range("C4").Value = 1.141827825028099
$ws.Range("D4").Value = 10.27012769846439
$ws.Range("E4").Value = 5.522519856084116
$ws.Range("F4").Value = 1.08577167714555
$ws.Range("G4").Value = 6.734269579350196
$ws.Range("H4").Value = 1.87267901678827
$ws.Range("I4").Value = 2.448790887616315
$ws.Range("J4").Value = 0.561130630218199
$ws.Range("K4").Value = 0.784136477480317
$ws.Range("L4").Value = 1.506840742146986
$ws.Range("M4").Value = 0.05218428679276705
$ws.Range("N4").Value = 0.4418758080696275
$ws.Range("O4").Value = 0.3324865172726089
$ws.Range("B5").Value = 67.7697631628386
$ws.Range("C5").Value = 1.15894054821536
$ws.Range("D5").Value = 10.28354262044705
$ws.Range("E5").Value = 5.540726374241397
$ws.Range("F5").Value = 1.090547181025713
$ws.Range("G5").Value = 6.733830960356887
$ws.Range("H5").Value = 1.889120606924848
$ws.Range("I5").Value = 2.449373899252271
$ws.Range("J5").Value = 0.5574774803444346
$ws.Range("K5").Value = 0.7991214305753612
$ws.Range("L5").Value = 1.507965383697544
$ws.Range("M5").Value = 0.05245194856038579
$ws.Range("N5").Value = 0.4221045815449641
$ws.Range("O5").Value = 0.3305283256744951
$ws.Range("B6").Value = 67.75925070470859
$ws.Range("C6").Value = 1.154598400614934
$ws.Range("D6").Value = 10.27754898620034
$ws.Range("E6").Value = 5.530274469708507
$ws.Range("F6").Value = 1.102364766951661
$ws.Range("G6").Value = 6.772207614083779
$ws.Range("H6").Value = 1.919869424809765
$ws.Range("I6").Value = 2.443909979980434
$ws.Range("J6").Value = 0.5624447249516811
$ws.Range("K6").Value = 0.7755157123271125
$ws.Range("L6").Value = 1.530534811517581
$ws.Range("M6").Value = 0.052635142840614
$ws.Range("N6").Value = 0.4259398547171418
$ws.Range("O6").Value = 0.3400575753830432
$ws.Range("B7").Value = 67.75472752030541
$ws.Range("C7").Value = 1.160609896621011
$ws.Range("D7").Value = 10.29634946761374
$ws.Range("E7").Value = 5.545704424467215
$ws.Range("F7").Value = 1.098277886123969
$ws.Range("G7").Value = 6.763775082207779
$ws.Range("H7").Value = 1.895476673740986
$ws.Range("I7").Value = 2.436803013971721
$ws.Range("J7").Value = 0.5601408894718432
$ws.Range("K7").Value = 0.8024353276403181
$ws.Range("L7").Value = 1.506438602539693
$ws.Range("M7").Value = 0.05339251125979835
$ws.Range("N7").Value = 0.4263757600638173
$ws.Range("O7").Value = 0.3292133670421802
